$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-09-07 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-08 Sunday", 2)
$d.Content.Find.Execute("346÷3=115, 1", $true, $false, $false, $false, $false, $true, 1, $false, "282÷8=35, 2", 2)
$d.Content.Find.Execute("891÷9=99, 0", $true, $false, $false, $false, $false, $true, 1, $false, "543÷7=77, 4", 2)
$d.Content.Find.Execute("254÷6=42, 2", $true, $false, $false, $false, $false, $true, 1, $false, "548÷2=274, 0", 2)
$d.Content.Find.Execute("980÷8=122, 4", $true, $false, $false, $false, $false, $true, 1, $false, "806÷2=403, 0", 2)
$d.Content.Find.Execute("642÷3=214, 0", $true, $false, $false, $false, $false, $true, 1, $false, "431÷2=215, 1", 2)
$d.Content.Find.Execute("208÷5=41, 3", $true, $false, $false, $false, $false, $true, 1, $false, "216÷6=36, 0", 2)
$d.Content.Find.Execute("241÷3=80, 1", $true, $false, $false, $false, $false, $true, 1, $false, "161÷4=40, 1", 2)
$d.Content.Find.Execute("611÷2=305, 1", $true, $false, $false, $false, $false, $true, 1, $false, "651÷9=72, 3", 2)
$d.Content.Find.Execute("947÷3=315, 2", $true, $false, $false, $false, $false, $true, 1, $false, "908÷8=113, 4", 2)
$d.Content.Find.Execute("979÷6=163, 1", $true, $false, $false, $false, $false, $true, 1, $false, "262÷4=65, 2", 2)
$d.Content.Find.Execute("801÷4=200, 1", $true, $false, $false, $false, $false, $true, 1, $false, "253÷8=31, 5", 2)
$d.Content.Find.Execute("746÷2=373, 0", $true, $false, $false, $false, $false, $true, 1, $false, "309÷7=44, 1", 2)
$d.Content.Find.Execute("295÷7=42, 1", $true, $false, $false, $false, $false, $true, 1, $false, "649÷2=324, 1", 2)
$d.Content.Find.Execute("916÷8=114, 4", $true, $false, $false, $false, $false, $true, 1, $false, "848÷2=424, 0", 2)
$d.Content.Find.Execute("355÷9=39, 4", $true, $false, $false, $false, $false, $true, 1, $false, "279÷9=31, 0", 2)
$d.Content.Find.Execute("748÷3=249, 1", $true, $false, $false, $false, $false, $true, 1, $false, "253÷7=36, 1", 2)
$d.Content.Find.Execute("126÷4=31, 2", $true, $false, $false, $false, $false, $true, 1, $false, "936÷6=156, 0", 2)
$d.Content.Find.Execute("336÷9=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "508÷7=72, 4", 2)
$d.Content.Find.Execute("479÷5=95, 4", $true, $false, $false, $false, $false, $true, 1, $false, "356÷7=50, 6", 2)
$d.Content.Find.Execute("847÷3=282, 1", $true, $false, $false, $false, $false, $true, 1, $false, "342÷8=42, 6", 2)
$d.Content.Find.Execute("529÷2=264, 1", $true, $false, $false, $false, $false, $true, 1, $false, "471÷5=94, 1", 2)
$d.Content.Find.Execute("554÷7=79, 1", $true, $false, $false, $false, $false, $true, 1, $false, "257÷8=32, 1", 2)
$d.Content.Find.Execute("296÷7=42, 2", $true, $false, $false, $false, $false, $true, 1, $false, "639÷4=159, 3", 2)
$d.Content.Find.Execute("668÷2=334, 0", $true, $false, $false, $false, $false, $true, 1, $false, "567÷5=113, 2", 2)
$d.Content.Find.Execute("177÷9=19, 6", $true, $false, $false, $false, $false, $true, 1, $false, "382÷8=47, 6", 2)
